$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.039.30"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.828.27"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "312.44"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "0.4605"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("D8").Value = "0.3698"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "0.07351"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").Value = "0.8716"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "0.07956"
$ws.Range("E11").Value = "  +3.35%  "
$ws.Range("D12").Value = "19.76"
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("D13").Value = "5.342"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "6.551"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "91.77"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").Value = "1.705.16"
$ws.Range("E16").Value = "  -8.14%  "
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008871"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "14.68"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "26.901.74"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").Value = "5.125"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").Value = "10.58"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "2.000.41"
$ws.Range("E24").Value = "  -4.29%  "
$ws.Range("D25").Value = "152.47"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "1.846"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("D27").Value = "18.53"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "2.071"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").Value = "5.088"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.30"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "0.08879"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "0.7315"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.440"
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "1.134"
$ws.Range("E35").Value = "  -2.67%  "
$ws.Range("D36").Value = "2.455"
$ws.Range("E36").Value = "  -3.01%  "
$ws.Range("D37").Value = "1.074"
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.05244"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01943"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "2.946"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "7.129"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D42").Value = "0.5163"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "0.8906"
$ws.Range("E43").Value = "  -11.94%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1630"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.210"
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("D46").Value = "0.4827"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "1.009"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "10.25"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").Value = "102.07"
$ws.Range("E49").Value = "  -2.46%  "
$ws.Range("D50").Value = "1.626"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("E51").Value = "  -0.74%  "
